# Adds crime-count data for 2025-12-10 by updating the 2025 year-to-date
# totals (column L) across the Citywide Totals, By Neighborhood, and all
# per-neighborhood worksheets.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{Sheet="Citywide Totals"; Cell="L2"; Value=6237},
    @{Sheet="Citywide Totals"; Cell="L3"; Value=6769},
    @{Sheet="Citywide Totals"; Cell="L4"; Value=1675},
    @{Sheet="Citywide Totals"; Cell="L5"; Value=400},
    @{Sheet="Citywide Totals"; Cell="L6"; Value=5560},
    @{Sheet="Citywide Totals"; Cell="L7"; Value=20641},
    @{Sheet="Logan Square"; Cell="L4"; Value=28},
    @{Sheet="Logan Square"; Cell="L7"; Value=232},
    @{Sheet="Austin"; Cell="L2"; Value=415},
    @{Sheet="Austin"; Cell="L3"; Value=484},
    @{Sheet="Austin"; Cell="L7"; Value=1364},
    @{Sheet="South Chicago"; Cell="L2"; Value=142},
    @{Sheet="South Chicago"; Cell="L6"; Value=101},
    @{Sheet="South Chicago"; Cell="L7"; Value=450},
    @{Sheet="Garfield Park"; Cell="L2"; Value=253},
    @{Sheet="Garfield Park"; Cell="L3"; Value=328},
    @{Sheet="Garfield Park"; Cell="L6"; Value=265},
    @{Sheet="Garfield Park"; Cell="L7"; Value=931},
    @{Sheet="Grand Crossing"; Cell="L3"; Value=282},
    @{Sheet="Grand Crossing"; Cell="L7"; Value=793},
    @{Sheet="New City"; Cell="L2"; Value=146},
    @{Sheet="New City"; Cell="L7"; Value=403},
    @{Sheet="By Neighborhood"; Cell="L2"; Value=183},
    @{Sheet="By Neighborhood"; Cell="L4"; Value=73},
    @{Sheet="By Neighborhood"; Cell="L6"; Value=165},
    @{Sheet="By Neighborhood"; Cell="L7"; Value=659},
    @{Sheet="By Neighborhood"; Cell="L8"; Value=1364},
    @{Sheet="By Neighborhood"; Cell="L18"; Value=141},
    @{Sheet="By Neighborhood"; Cell="L19"; Value=558},
    @{Sheet="By Neighborhood"; Cell="L22"; Value=68},
    @{Sheet="By Neighborhood"; Cell="L29"; Value=1151},
    @{Sheet="By Neighborhood"; Cell="L33"; Value=931},
    @{Sheet="By Neighborhood"; Cell="L37"; Value=793},
    @{Sheet="By Neighborhood"; Cell="L48"; Value=271},
    @{Sheet="By Neighborhood"; Cell="L49"; Value=112},
    @{Sheet="By Neighborhood"; Cell="L50"; Value=101},
    @{Sheet="By Neighborhood"; Cell="L53"; Value=232},
    @{Sheet="By Neighborhood"; Cell="L55"; Value=220},
    @{Sheet="By Neighborhood"; Cell="L58"; Value=11},
    @{Sheet="By Neighborhood"; Cell="L60"; Value=137},
    @{Sheet="By Neighborhood"; Cell="L63"; Value=66},
    @{Sheet="By Neighborhood"; Cell="L65"; Value=403},
    @{Sheet="By Neighborhood"; Cell="L67"; Value=717},
    @{Sheet="By Neighborhood"; Cell="L75"; Value=76},
    @{Sheet="By Neighborhood"; Cell="L76"; Value=320},
    @{Sheet="By Neighborhood"; Cell="L78"; Value=269},
    @{Sheet="By Neighborhood"; Cell="L79"; Value=570},
    @{Sheet="By Neighborhood"; Cell="L83"; Value=450},
    @{Sheet="By Neighborhood"; Cell="L84"; Value=198},
    @{Sheet="By Neighborhood"; Cell="L85"; Value=1026},
    @{Sheet="By Neighborhood"; Cell="L88"; Value=220},
    @{Sheet="By Neighborhood"; Cell="L91"; Value=278},
    @{Sheet="By Neighborhood"; Cell="L93"; Value=105},
    @{Sheet="By Neighborhood"; Cell="L94"; Value=252},
    @{Sheet="By Neighborhood"; Cell="L96"; Value=229},
    @{Sheet="By Neighborhood"; Cell="L101"; Value=20641},
    @{Sheet="North Lawndale"; Cell="L5"; Value=21},
    @{Sheet="North Lawndale"; Cell="L6"; Value=166},
    @{Sheet="North Lawndale"; Cell="L7"; Value=717},
    @{Sheet="South Deering"; Cell="L2"; Value=66},
    @{Sheet="South Deering"; Cell="L7"; Value=198},
    @{Sheet="Lincoln Park"; Cell="L6"; Value=45},
    @{Sheet="Lincoln Park"; Cell="L7"; Value=112},
    @{Sheet="Englewood"; Cell="L3"; Value=444},
    @{Sheet="Englewood"; Cell="L7"; Value=1151},
    @{Sheet="Lake View"; Cell="L4"; Value=52},
    @{Sheet="Lake View"; Cell="L7"; Value=271},
    @{Sheet="Chatham"; Cell="L2"; Value=202},
    @{Sheet="Chatham"; Cell="L3"; Value=171},
    @{Sheet="Chatham"; Cell="L7"; Value=558},
    @{Sheet="River North"; Cell="L4"; Value=38},
    @{Sheet="River North"; Cell="L7"; Value=320},
    @{Sheet="Ashburn"; Cell="L3"; Value=49},
    @{Sheet="Ashburn"; Cell="L7"; Value=165},
    @{Sheet="Rogers Park"; Cell="L3"; Value=89},
    @{Sheet="Rogers Park"; Cell="L7"; Value=269},
    @{Sheet="Lower West Side"; Cell="L2"; Value=65},
    @{Sheet="Lower West Side"; Cell="L7"; Value=220},
    @{Sheet="West Ridge"; Cell="L3"; Value=66},
    @{Sheet="West Ridge"; Cell="L7"; Value=229},
    @{Sheet="Washington Park"; Cell="L3"; Value=126},
    @{Sheet="Washington Park"; Cell="L7"; Value=278},
    @{Sheet="Roseland"; Cell="L2"; Value=178},
    @{Sheet="Roseland"; Cell="L7"; Value=570},
    @{Sheet="Calumet Heights"; Cell="L3"; Value=50},
    @{Sheet="Calumet Heights"; Cell="L7"; Value=141},
    @{Sheet="West Lawn"; Cell="L2"; Value=36},
    @{Sheet="West Lawn"; Cell="L7"; Value=105},
    @{Sheet="Auburn Gresham"; Cell="L3"; Value=207},
    @{Sheet="Auburn Gresham"; Cell="L7"; Value=659},
    @{Sheet="West Loop"; Cell="L6"; Value=93},
    @{Sheet="West Loop"; Cell="L7"; Value=252},
    @{Sheet="Lincoln Square"; Cell="L2"; Value=33},
    @{Sheet="Lincoln Square"; Cell="L7"; Value=101},
    @{Sheet="Albany Park"; Cell="L2"; Value=60},
    @{Sheet="Albany Park"; Cell="L7"; Value=183},
    @{Sheet="United Center"; Cell="L6"; Value=60},
    @{Sheet="United Center"; Cell="L7"; Value=220},
    @{Sheet="Pullman"; Cell="L3"; Value=26},
    @{Sheet="Pullman"; Cell="L7"; Value=76},
    @{Sheet="Morgan Park"; Cell="L2"; Value=46},
    @{Sheet="Morgan Park"; Cell="L3"; Value=45},
    @{Sheet="Morgan Park"; Cell="L7"; Value=137},
    @{Sheet="South Shore"; Cell="L3"; Value=427},
    @{Sheet="South Shore"; Cell="L7"; Value=1026},
    @{Sheet="Clearing"; Cell="L3"; Value=27},
    @{Sheet="Clearing"; Cell="L7"; Value=68},
    @{Sheet="Archer Heights"; Cell="L6"; Value=22},
    @{Sheet="Archer Heights"; Cell="L7"; Value=73},
    @{Sheet="Millenium Park"; Cell="L6"; Value=5},
    @{Sheet="Millenium Park"; Cell="L7"; Value=11}
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}

Write-Output "Applied $($updates.Count) cell updates."
